{"js": "// Revert \"Modiied FlutterTutorial Doc\":\n// Remove the trailing \"two compilation processes\" discussion (the Yes/Dart\n// Compilation/Flutter Compilation/So in summary paragraphs, plus their\n// separating blank paragraphs) that followed the \"By having only one\n// codebase...\" paragraph, restoring that paragraph as the last one in the\n// document body.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\n// Text that marks the end of the content we want to KEEP. Everything after\n// this paragraph (until the end of the body) should be removed.\nconst keepUntilMarker = \"By having only one codebase\";\n// Text that marks the last paragraph we want to REMOVE (inclusive).\nconst lastRemoveMarker = \"So, in summary, both Dart and Flutter\";\n\nlet keepIndex = -1;\nlet lastRemoveIndexFound = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const text = paragraphs.items[i].text;\n  if (keepIndex === -1 && text.indexOf(keepUntilMarker) !== -1) {\n    keepIndex = i;\n  }\n  if (text.indexOf(lastRemoveMarker) !== -1) {\n    lastRemoveIndexFound = i;\n  }\n}\n\nif (keepIndex !== -1 && lastRemoveIndexFound !== -1 && lastRemoveIndexFound > keepIndex) {\n  const toDelete = [];\n  for (let i = keepIndex + 1; i <= lastRemoveIndexFound; i++) {\n    toDelete.push(paragraphs.items[i]);\n  }\n  for (const para of toDelete) {\n    para.delete();\n  }\n  await context.sync();\n}\n", "ps1": "# Revert \"Modiied FlutterTutorial Doc\":\n# Remove the trailing \"two compilation processes\" discussion (the Yes/Dart\n# Compilation/Flutter Compilation/So in summary paragraphs, plus their\n# separating blank paragraphs) that followed the \"By having only one\n# codebase...\" paragraph, restoring that paragraph as the last one in the\n# document body.\n\n$d = $word.ActiveDocument\n\n$keepUntilMarker = \"By having only one codebase\"\n$lastRemoveMarker = \"So, in summary, both Dart and Flutter\"\n\n$count = $d.Paragraphs.Count\n\n$keepIndex = -1\n$lastRemoveIndex = -1\n\nfor ($i = 1; $i -le $count; $i++) {\n    $text = $d.Paragraphs.Item($i).Range.Text\n    if ($keepIndex -eq -1 -and $text -like \"*$keepUntilMarker*\") {\n        $keepIndex = $i\n    }\n    if ($text -like \"*$lastRemoveMarker*\") {\n        $lastRemoveIndex = $i\n    }\n}\n\nif ($keepIndex -ne -1 -and $lastRemoveIndex -ne -1 -and $lastRemoveIndex -gt $keepIndex) {\n    # Delete the whole block in a single Range.Delete() call (rather than\n    # paragraph-by-paragraph) so the very last paragraph mark of the\n    # document is correctly merged away instead of left behind empty.\n    $startRange = $d.Paragraphs.Item($keepIndex + 1).Range\n    $endRange = $d.Paragraphs.Item($lastRemoveIndex).Range\n    $rng = $d.Range($startRange.Start, $endRange.End)\n    $rng.Delete()\n}\n"}
